$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at position 91, shifting existing rows 91-123 down to 95-127
$ws.Range("A91:T94").EntireRow.Insert()

# Row 91: new price record, date 44726, variedad "Granny Smith"
$ws.Range("A91").Value = 1
$ws.Range("B91").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C91").Value = "Arica y Parinacota"
$ws.Range("D91").Value = 44726
$ws.Range("E91").Value = 15
$ws.Range("F91").Value = "Fruta"
$ws.Range("G91").Value = 100104
$ws.Range("H91").Value = "Frutos de pepita"
$ws.Range("I91").Value = 100104002
$ws.Range("J91").Value = "Manzana"
$ws.Range("K91").Value = "Granny Smith"
$ws.Range("L91").Value = "Calibre 90"
$ws.Range("M91").Value = 300
$ws.Range("N91").Value = 17000
$ws.Range("O91").Value = 18000
$ws.Range("P91").Value = 17500
$ws.Range("Q91").Value = "`$/caja 18 kilos embalada"
$ws.Range("R91").Value = "Región de O'Higgins"
$ws.Range("S91").Value = 972
$ws.Range("T91").Value = 18

# Row 92: new price record, date 44726, variedad "Pink Lady"
$ws.Range("A92").Value = 1
$ws.Range("B92").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C92").Value = "Arica y Parinacota"
$ws.Range("D92").Value = 44726
$ws.Range("E92").Value = 15
$ws.Range("F92").Value = "Fruta"
$ws.Range("G92").Value = 100104
$ws.Range("H92").Value = "Frutos de pepita"
$ws.Range("I92").Value = 100104002
$ws.Range("J92").Value = "Manzana"
$ws.Range("K92").Value = "Pink Lady"
$ws.Range("L92").Value = "Calibre 90"
$ws.Range("M92").Value = 270
$ws.Range("N92").Value = 17000
$ws.Range("O92").Value = 18000
$ws.Range("P92").Value = 17500
$ws.Range("Q92").Value = "`$/caja 18 kilos embalada"
$ws.Range("R92").Value = "Región de O'Higgins"
$ws.Range("S92").Value = 972
$ws.Range("T92").Value = 18

# Row 93: new price record, date 44726, variedad "Royal Gala"
$ws.Range("A93").Value = 1
$ws.Range("B93").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C93").Value = "Arica y Parinacota"
$ws.Range("D93").Value = 44726
$ws.Range("E93").Value = 15
$ws.Range("F93").Value = "Fruta"
$ws.Range("G93").Value = 100104
$ws.Range("H93").Value = "Frutos de pepita"
$ws.Range("I93").Value = 100104002
$ws.Range("J93").Value = "Manzana"
$ws.Range("K93").Value = "Royal Gala"
$ws.Range("L93").Value = "Calibre 80"
$ws.Range("M93").Value = 300
$ws.Range("N93").Value = 17000
$ws.Range("O93").Value = 18000
$ws.Range("P93").Value = 17500
$ws.Range("Q93").Value = "`$/caja 18 kilos embalada"
$ws.Range("R93").Value = "Región de O'Higgins"
$ws.Range("S93").Value = 972
$ws.Range("T93").Value = 18

# Row 94: new price record, date 44726, variedad "Scarlett"
$ws.Range("A94").Value = 1
$ws.Range("B94").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C94").Value = "Arica y Parinacota"
$ws.Range("D94").Value = 44726
$ws.Range("E94").Value = 15
$ws.Range("F94").Value = "Fruta"
$ws.Range("G94").Value = 100104
$ws.Range("H94").Value = "Frutos de pepita"
$ws.Range("I94").Value = 100104002
$ws.Range("J94").Value = "Manzana"
$ws.Range("K94").Value = "Scarlett"
$ws.Range("L94").Value = "Calibre 80"
$ws.Range("M94").Value = 300
$ws.Range("N94").Value = 17000
$ws.Range("O94").Value = 18000
$ws.Range("P94").Value = 17500
$ws.Range("Q94").Value = "`$/caja 18 kilos embalada"
$ws.Range("R94").Value = "Región de O'Higgins"
$ws.Range("S94").Value = 972
$ws.Range("T94").Value = 18

